$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking strings
# (e.g. "1.00", "29.60", "64.203.14") are not auto-coerced to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Price (D) updates
$ws.Range('D2').Value = '64.203.14'
$ws.Range('D3').Value = '2.524.11'
$ws.Range('D5').Value = '581.47'
$ws.Range('D6').Value = '152.52'
$ws.Range('D7').Value = '1.00'
$ws.Range('D8').Value = '0.537'
$ws.Range('D12').Value = '0.356'
$ws.Range('D13').Value = '29.60'
$ws.Range('D14').Value = '0.0000179'
$ws.Range('D15').Value = '2.982.58'
$ws.Range('D16').Value = '63.565.58'
$ws.Range('D17').Value = '2.520.51'
$ws.Range('D18').Value = '7.86'
$ws.Range('D19').Value = '10.98'
$ws.Range('D20').Value = '4.26'
$ws.Range('D21').Value = '327.87'
$ws.Range('D22').Value = '2.26'
$ws.Range('D25').Value = '65.50'
$ws.Range('D26').Value = '665.26'
$ws.Range('D29').Value = '1.49'
$ws.Range('D30').Value = '0.992'
$ws.Range('D32').Value = '1.86'
$ws.Range('D33').Value = '0.135'
$ws.Range('D36').Value = '4.82'
$ws.Range('D37').Value = '5.55'
$ws.Range('D38').Value = '0.372'
$ws.Range('D39').Value = '18.84'
$ws.Range('D40').Value = '151.97'
$ws.Range('D41').Value = '2.80'
$ws.Range('D43').Value = '0.999'
$ws.Range('D44').Value = '158.57'
$ws.Range('D45').Value = '0.0₆0301'
$ws.Range('D48').Value = '21.06'
$ws.Range('D49').Value = '0.618'
$ws.Range('D50').Value = '0.0521'
$ws.Range('D51').Value = '0.0229'

# Restore original (unset) number formatting now that values are stored as text
$ws.Range("D2:D51").ClearFormats()

# Coin / Link / Volume(1h) updates
$ws.Range('E2').Value = '  +2.19%  '
$ws.Range('E3').Value = '  +2.50%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('E5').Value = '  +1.20%  '
$ws.Range('E6').Value = '  +4.66%  '
$ws.Range('E7').Value = '  +0.01%  '
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  +1.26%  '
$ws.Range('E10').Value = '  +0.08%  '
$ws.Range('E11').Value = '  +0.29%  '
$ws.Range('E12').Value = '  +0.08%  '
$ws.Range('E13').Value = '  +2.35%  '
$ws.Range('E14').Value = '  +1.57%  '
$ws.Range('E15').Value = '  +2.48%  '
$ws.Range('E16').Value = '  +1.33%  '
$ws.Range('E17').Value = '  +2.29%  '
$ws.Range('E18').Value = '  -2.24%  '
$ws.Range('E19').Value = '  +0.09%  '
$ws.Range('E20').Value = '  +3.02%  '
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('E22').Value = '  +1.93%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('E24').Value = '  -0.94%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('E25').Value = '  -0.26%  '
$ws.Range('B26').Value = 'Bittensor'
$ws.Range('C26').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('E26').Value = '  +1.56%  '
$ws.Range('E27').Value = '  +3.30%  '
$ws.Range('E29').Value = '  +2.47%  '
$ws.Range('E30').Value = '  -0.89%  '
$ws.Range('E31').Value = '  +0.60%  '
$ws.Range('E32').Value = '  +0.65%  '
$ws.Range('E33').Value = '  +0.74%  '
$ws.Range('E34').Value = '  -0.18%  '
$ws.Range('E35').Value = '  +0.93%  '
$ws.Range('E36').Value = '  +1.39%  '
$ws.Range('E37').Value = '  +3.46%  '
$ws.Range('E38').Value = '  +0.58%  '
$ws.Range('E39').Value = '  +0.80%  '
$ws.Range('B40').Value = 'Monero'
$ws.Range('C40').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('E40').Value = '  +0.71%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('E41').Value = '  +2.28%  '
$ws.Range('E42').Value = '  +3.22%  '
$ws.Range('E43').Value = '  +0.01%  '
$ws.Range('E44').Value = '  +2.82%  '
$ws.Range('E45').Value = '  -3.86%  '
$ws.Range('E46').Value = '  +1.36%  '
$ws.Range('E47').Value = '  +1.21%  '
$ws.Range('E48').Value = '  +3.75%  '
$ws.Range('E49').Value = '  +1.91%  '
$ws.Range('E50').Value = '  +2.03%  '
$ws.Range('E51').Value = '  +1.85%  '
